$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: new log entry (4794 / 01/04/2020 / times / help debug) ---
$ws.Cells.Item(10, 2).Value = 4794

$fmt = $ws.Range("C10").NumberFormat
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "01/04/2020"
$ws.Range("C10").NumberFormat = $fmt

$ws.Range("D10").Value = 0.60416666666666663
$ws.Range("E10").Value = 0.66666666666666663
$ws.Range("G10").Value = "Help debug arithmetic unit"

# --- Row 9: update description to reflect work done that day ---
$ws.Range("G9").Value = "worked Arithmetic unit"

# --- Row 11: new log entry (02/04/2020) ---
$fmt = $ws.Range("C11").NumberFormat
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "02/04/2020"
$ws.Range("C11").NumberFormat = $fmt

$ws.Range("D11").Value = 0.72916666666666663
$ws.Range("E11").Value = 0.77083333333333337
$ws.Range("G11").Value = "Screenshots of waves of functional simulation"

# --- Row 12: another entry on 02/04/2020 ---
$fmt = $ws.Range("C12").NumberFormat
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "02/04/2020"
$ws.Range("C12").NumberFormat = $fmt

$ws.Range("D12").Value = 0.81944444444444453
$ws.Range("E12").Value = 0.94791666666666663
$ws.Range("G12").Value = "Worked on documentations"

# --- Update the view: scroll so column C is left-most, select C13 ---
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("C13").Select()
